$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '89.685.12'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +4.16%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.231.50'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.90%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.20'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.62%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '620.99'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.30%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.405'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +12.27%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.707'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +9.33%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.09%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.224.27'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.98%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.576'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.57%  '

# Row 12
$ws.Range('B12').Value = 'ShibaInu'
$ws.Range('C12').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000275'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +8.57%  '

# Row 13
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.179'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.93%  '

# Row 14
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.57%  '

# Row 15
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '33.60'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.49%  '

# Row 16
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.823.17'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.77%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.447.25'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.03%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.248.38'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.03%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.21'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +8.41%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.93'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.47%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '423.86'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.61%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.81'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.16%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.21'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.80%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000182'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +42.11%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.37'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +6.07%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.95'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.10%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.379.01'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.33%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '75.23'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.20%  '

# Row 29
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.07%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.170'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.50%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.27%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '559.87'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.17%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '8.50'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.86%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.08'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.75%  '

# Row 35
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.11%  '

# Row 36
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.30%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.50'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +20.30%  '

# Row 38
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.133'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.08%  '

# Row 39
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '22.29'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.06%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '21.91'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.46%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.08%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.388'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.20%  '

# Row 43
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.98'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.58%  '

# Row 44
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.12%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '149.87'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.22%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '180.56'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.20%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '43.89'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.01%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.128'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +8.81%  '

# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.40%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.617'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.05%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.77'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.12%  '
